# Append the newly-recorded "Spezialpreise" (special price) sales entries for
# the Tabubrecherin event and the Weihnachten (Christmas) market in der
# Schustergasse to the Table1 listing on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# -- Row 31: 08.12.2024 - Spez 1 / Kaffee und Gipfeli --------------------
$lo.ListRows.Add() | Out-Null
$ws.Range("A31").Value = 45634
$ws.Range("B31").Value = "Spez 1"
$ws.Range("C31").Value = "Kaffee und Gipfeli"
$ws.Range("D31").Value = 5

# -- Row 32: 15.12.2024 - Spez 1 / Glühwein -------------------------------
$lo.ListRows.Add() | Out-Null
$ws.Range("A32").Value = 45641
$ws.Range("B32").Value = "Spez 1"
$ws.Range("C32").Value = "Glühwein"
$ws.Range("D32").Value = 3

# -- Row 33: 15.12.2024 - Spez 2 / Punsch ---------------------------------
$lo.ListRows.Add() | Out-Null
$ws.Range("A33").Value = 45641
$ws.Range("B33").Value = "Spez 2"
$ws.Range("C33").Value = "Punsch"
$ws.Range("D33").Value = 2

# -- Row 34: 15.12.2024 - Spez 3 / Tee ------------------------------------
$lo.ListRows.Add() | Out-Null
$ws.Range("A34").Value = 45641
$ws.Range("B34").Value = "Spez 3"
$ws.Range("C34").Value = "Tee"
$ws.Range("D34").Value = 2

# -- Row 35: 15.12.2024 - Spez 4 / Wienerli und Brot ----------------------
$lo.ListRows.Add() | Out-Null
$ws.Range("A35").Value = 45641
$ws.Range("B35").Value = "Spez 4"
$ws.Range("C35").Value = "Wienerli und Brot"
$ws.Range("D35").Value = 4

# -- Row 36: 15.12.2024 - Spez 5 / Raclette - Schnitte --------------------
$lo.ListRows.Add() | Out-Null
$ws.Range("A36").Value = 45641
$ws.Range("B36").Value = "Spez 5"
$ws.Range("C36").Value = "Raclette - Schnitte"
$ws.Range("D36").Value = 4

# Move view/selection down to the next empty row, same as the author left it.
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A37").Select() | Out-Null
